$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width changes ---
# Stored OOXML width = ((requested*MDW+5)/MDW) rounded; MDW=6 here, so to
# land exactly on an integer stored width we request (target - 5/6).
$ws.Columns.Item(2).ColumnWidth = 46 - (5/6)
$ws.Columns.Item(5).ColumnWidth = 16 - (5/6)

# --- New row 5 ---
$ws.Range("A5").Value = "f4 address"
$ws.Range("B5").Value = "f4 city"
$ws.Range("C5").Value = "ff4 first name"
$ws.Range("D5").Value = "f4 last name"

$ws.Range("E5").Value = "'0404959583938"
$ws.Range("E5").Style = "Normal"

$ws.Range("F5").Value = " "

$ws.Range("G5").Value = "'040404"
$ws.Range("G5").Style = "Normal"

# --- New row 6 ---
$ws.Range("A6").Value = "f2 address"
$ws.Range("B6").Value = 'f2 Ho Chi MInh )(*#^@%@&@15`~=---:,.?"|\n \t'
$ws.Range("C6").Value = "f2 first"
$ws.Range("D6").Value = "f2 last"

$ws.Range("E6").Value = "'02738483925363"
$ws.Range("E6").Style = "Normal"

$ws.Range("F6").Value = "'0202"
$ws.Range("F6").Style = "Normal"

$ws.Range("G6").Value = "'020202"
$ws.Range("G6").Style = "Normal"
